$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new value "waltuh" in A2 (new shared string, referenced from A2)
$ws.Range("A2").Value = "waltuh"

# Update the selection to F5 (matches selection diff in sheetView)
$ws.Range("F5").Select()
